# Fills in status_done, status_checked, analyzed_visual (and node / date
# submitted where missing) for the newly-finished sample rows (20-55),
# mirroring the "very large and late update" commit: preprocessing runs
# that finished on 21.10.2022 / 23.10.2022 / 24.10.2022 are marked done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "yes"
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = "yes"
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = "yes"
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = "yes "
$ws.Range("D31").Value = 6
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = "yes"
$ws.Range("D32").Value = 6
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = "yes"
$ws.Range("D33").Value = 6
$ws.Range("E33").Value = 2
$ws.Range("F33").Value = "yes "
$ws.Range("D34").Value = 7
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = "yes"
$ws.Range("D35").Value = 6
$ws.Range("E35").Value = 2
$ws.Range("F35").Value = "yes"
$ws.Range("D36").Value = 6
$ws.Range("E36").Value = 2
$ws.Range("F36").Value = "yes"
$ws.Range("D37").Value = 6
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = "yes"
$ws.Range("D38").Value = 6
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = "yes"
$ws.Range("D39").Value = 6
$ws.Range("E39").Value = 2
$ws.Range("F39").Value = "yes"
$ws.Range("D40").Value = 6
$ws.Range("E40").Value = 2
$ws.Range("F40").Value = "yes"
$ws.Range("D41").Value = 6
$ws.Range("E41").Value = 2
$ws.Range("F41").Value = "yes"
$ws.Range("D42").Value = 6
$ws.Range("E42").Value = 2
$ws.Range("F42").Value = "yes"
$ws.Range("D43").Value = 6
$ws.Range("E43").Value = 2
$ws.Range("F43").Value = "yes "
$ws.Range("D44").Value = 6
$ws.Range("E44").Value = 2
$ws.Range("F44").Value = "yes"
$ws.Range("D45").Value = 6
$ws.Range("E45").Value = 2
$ws.Range("F45").Value = "yes"
$ws.Range("D46").Value = 6
$ws.Range("E46").Value = 2
$ws.Range("F46").Value = "yes"
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = "21.10.2022"
$ws.Range("D48").Value = 6
$ws.Range("E48").Value = 2
$ws.Range("F48").Value = "yes"
$ws.Range("B49").Value = 3
$ws.Range("C49").Value = "21.10.2022"
$ws.Range("D49").Value = 6
$ws.Range("E49").Value = 2
$ws.Range("F49").Value = "yes"
$ws.Range("B50").Value = 4
$ws.Range("C50").Value = "21.10.2022"
$ws.Range("D50").Value = 6
$ws.Range("E50").Value = 2
$ws.Range("F50").Value = "yes"
$ws.Range("B51").Value = 5
$ws.Range("C51").Value = "21.10.2022"
$ws.Range("D51").Value = 6
$ws.Range("E51").Value = 2
$ws.Range("F51").Value = "yes"
$ws.Range("B47").Value = 1
$ws.Range("C47").Value = "23.10.2022"
$ws.Range("D47").Value = 6
$ws.Range("E47").Value = 2
$ws.Range("F47").Value = "yes"
$ws.Range("B52").Value = 2
$ws.Range("C52").Value = "23.10.2022"
$ws.Range("D52").Value = 6
$ws.Range("E52").Value = 2
$ws.Range("F52").Value = "yes"
$ws.Range("B53").Value = 4
$ws.Range("C53").Value = "23.10.2022"
$ws.Range("D53").Value = 6
$ws.Range("E53").Value = 2
$ws.Range("F53").Value = "yes"
$ws.Range("B54").Value = 3
$ws.Range("C54").Value = "24.10.2022"
$ws.Range("D54").Value = 6
$ws.Range("E54").Value = 2
$ws.Range("F54").Value = "yes"
$ws.Range("B55").Value = 5
$ws.Range("C55").Value = "24.10.2022"
$ws.Range("D55").Value = 6
$ws.Range("E55").Value = 2
$ws.Range("F55").Value = "yes"

# Move the viewport/selection to where the user last left off editing.
$ws.Range("G52").Select()
